$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 88, pushing the existing rows 88..131 down to 89..132.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with this week's new record.
# Columns A,B,C,E,F,G,I,N,Q,R are constant across every data row in this sheet,
# so reuse them for the new row as well.
$ws.Range("A88").Value = 7
$ws.Range("B88").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C88").Value = "Ñuble"
$ws.Range("D88").Value = 44489
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112017
$ws.Range("G88").Value = "Apio"
$ws.Range("H88").Value = "Americana (o)"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 180
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 8500
$ws.Range("M88").Value = 8250
$ws.Range("N88").Value = "`$/docena de matas"
$ws.Range("O88").Value = "Provincia del Elquí"
$ws.Range("P88").Value = 1375
$ws.Range("Q88").Value = 6
$ws.Range("R88").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# every other date cell in column D.
$ws.Range("D88").NumberFormat = $ws.Range("D89").NumberFormat
